# "Agregado de estado de tarea"
# Sets the "Estatus" (column F) for every task row on the "Casos de Uso"
# sheet. Every row goes from "Por iniciar" to "Terminado", except the
# rows that correspond to the "Diagrama de secuencia" / odd task-steps
# (22, 23, 42, 43, 52) which are marked "No aplica" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Rows that should be flagged "No aplica" instead of "Terminado".
$noAplicaRows = @(22, 23, 42, 43, 52)

for ($row = 6; $row -le 56; $row++) {
    if ($noAplicaRows -contains $row) {
        $ws.Cells.Item($row, 6).Value = "No aplica"
    } else {
        $ws.Cells.Item($row, 6).Value = "Terminado"
    }
}

# Restore the selection on the visible sheet to match the saved view.
$ws.Activate()
$ws.Range("F52").Select()
